# Nuevo formato 15 jun 2021
# Update statistics and remove a student who is no longer "rescatable".

$wb = $excel.ActiveWorkbook

# --- Estadisticos 2P: update stats for TOMA MUESTRAS BIOLOGICAS / 2ALCM (row 3) ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 5
$ws2.Range("F3").Value = 36
$ws2.Range("G3").Value = 87.8

# --- Estadisticos Final: update stats for all three groups ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

# Row 2: PREPARA SOLUCIONES... / 2ALCM
$ws3.Range("E2").Value = 3
$ws3.Range("F2").Value = 38
$ws3.Range("G2").Value = 92.68
$ws3.Range("H2").Value = 7.3

# Row 3: TOMA MUESTRAS BIOLOGICAS / 2ALCM
$ws3.Range("E3").Value = 3
$ws3.Range("F3").Value = 38
$ws3.Range("G3").Value = 92.68
$ws3.Range("H3").Value = 7.5

# Row 4: TOMA MUESTRAS BIOLOGICAS / 2BLCM
$ws3.Range("E4").Value = 9
$ws3.Range("F4").Value = 27
$ws3.Range("G4").Value = 75
$ws3.Range("H4").Value = 6.9

# --- Rescatables: remove the row for student 20330051920246 (RAMOS RAMOS ISABELLA), ---
# --- who is no longer in need of a makeup; remaining rows shift up.                 ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows(3).Delete()
